$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") for data rows 2 through 120 from 45206 to 45208
for ($row = 2; $row -le 120; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45206) {
        $cell.Value2 = 45208
    }
}
